$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final ordered list of word-cloud terms (column A) with their weights (column B),
# after fixing typos ("sckit-learn" -> "scikit-learn", "LightBGM" -> "LightGBM"),
# replacing the duplicate "Gradient Boost" entry with "GBDT",
# normalizing case ("Web design" -> "Web Design", "Web scrapping" -> "Web Scrapping"),
# and re-ordering a few rows.
$words = @(
    "Python",
    "scikit-learn",
    "C++",
    "R",
    "Data Science",
    "Data Analytics",
    "Data Engineering",
    "ETL",
    "EDA",
    "Random Forest",
    "Machine Learning",
    "Deep Learning",
    "NLP",
    "Image Processing",
    "CNN",
    "TensorFlow",
    "Docker",
    "Git",
    "AWS",
    "SVM",
    "GBDT",
    "PCA",
    "XGBOOST",
    "LightGBM",
    "NLTK",
    "Polyglot",
    "Pandas",
    "Web/Data API",
    "MySQL",
    "MongoDB",
    "Web Design",
    "GUI",
    "html/css/javascript",
    "Keras",
    "Web Scrapping"
)

$weights = @(100,80,50,20,85,83,82,60,60,68,90,80,60,30,40,42,78,70,82,68,75,77,40,70,50,40,70,80,82,30,60,40,40,50,70)

# Write the brand-new terms first (in this exact order) so that they are
# appended to the shared-string table in the same order as in the target
# workbook: scikit-learn, Web Design, Web Scrapping, LightGBM, GBDT.
$ws.Cells.Item(2, 1).Value = "scikit-learn"
$ws.Cells.Item(31, 1).Value = "Web Design"
$ws.Cells.Item(35, 1).Value = "Web Scrapping"
$ws.Cells.Item(24, 1).Value = "LightGBM"
$ws.Cells.Item(21, 1).Value = "GBDT"

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $words[$i]
    $ws.Cells.Item($row, 2).Value = $weights[$i]
}

# Match the new column width and selection from the diff.
# (The engine's stored column width is ColumnWidth + 5/6, so back that out
# here to land on exactly "29" in the saved XML.)
$ws.Columns.Item(1).ColumnWidth = 28.1666666666667
$ws.Range("A4:B35").Select()
